# Apply the edits described by the commit: "Added bar codes and recovery dates where necessary"
$wb = $excel.ActiveWorkbook
$moorings = $wb.Worksheets.Item("Moorings")
$assetCal = $wb.Worksheets.Item("Asset_Cal_Info")

# --- Moorings sheet: fill in the Recover Date for the single deployment row ---
$moorings.Range("G2").Value = 42510

# --- Asset_Cal_Info sheet ---
# A lookup list of the unique Ref Des values was added in column L (rows 2-10),
# used as a cross-check helper against column A via MATCH() formulas in columns K and P.
$assetCal.Range("L2").Value = "CP02PMCI-SBS01-00-STCENG000"
$assetCal.Range("L3").Value = "CP02PMCI-SBS01-01-MOPAK0000"
$assetCal.Range("L4").Value = "CP02PMCI-RII01-02-ADCPTG010"
$assetCal.Range("L5").Value = "CP02PMCI-WFP01-00-WFPENG000"
$assetCal.Range("L6").Value = "CP02PMCI-WFP01-01-VEL3DK000"
$assetCal.Range("L7").Value = "CP02PMCI-WFP01-02-DOFSTK000"
$assetCal.Range("L8").Value = "CP02PMCI-WFP01-03-CTDPFK000"
$assetCal.Range("L9").Value = "CP02PMCI-WFP01-04-FLORTK000"
$assetCal.Range("L10").Value = "CP02PMCI-WFP01-05-PARADK000"

# K2:K39 - for every Ref Des row, find its position in the L lookup list
for ($r = 2; $r -le 39; $r++) {
    $assetCal.Range("K$r").Formula = "=MATCH(A$r,L:L,0)"
}

# P2:P10 - for every L lookup entry, find where it appears in column A
for ($r = 2; $r -le 10; $r++) {
    $assetCal.Range("P$r").Formula = "=MATCH(L$r,A:A,0)"
}

# The missing sensor barcode for the STCENG000 row was filled in
$assetCal.Range("E39").Value = "OL000341"
$assetCal.Range("E39").Style = "Normal"

# --- Restore selection / active-sheet state to match where the edits left off ---
$moorings.Range("E11").Select()
$assetCal.Select()
$assetCal.Range("F47").Select()
